$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/5/2024  Through  8/11/2024"

# --- Precinct crime-stat numeric updates (rows 15-28) ---
$ws.Range("M15").Value = -71.428571428571
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -36.842105263157
$ws.Range("I16").Value = 99
$ws.Range("J16").Value = 116
$ws.Range("K16").Value = -14.655172413793
$ws.Range("L16").Value = -31.25
$ws.Range("M16").Value = 17.857142857142
$ws.Range("N16").Value = -80.776699029126
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -27.777777777777
$ws.Range("I17").Value = 102
$ws.Range("J17").Value = 135
$ws.Range("K17").Value = -24.444444444444
$ws.Range("L17").Value = -3.773584905660
$ws.Range("M17").Value = 78.947368421052
$ws.Range("N17").Value = -43.333333333333
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 14
$ws.Range("H18").Value = -36.363636363636
$ws.Range("I18").Value = 155
$ws.Range("J18").Value = 186
$ws.Range("K18").Value = -16.666666666666
$ws.Range("L18").Value = -43.636363636363
$ws.Range("M18").Value = 31.355932203389
$ws.Range("N18").Value = -69.246031746031
$ws.Range("C19").Value = 31
$ws.Range("D19").Value = 33
$ws.Range("E19").Value = -6.060606060606
$ws.Range("F19").Value = 92
$ws.Range("G19").Value = 103
$ws.Range("H19").Value = -10.679611650485
$ws.Range("I19").Value = 657
$ws.Range("J19").Value = 762
$ws.Range("K19").Value = -13.779527559055
$ws.Range("L19").Value = -17.669172932330
$ws.Range("M19").Value = 0.921658986175
$ws.Range("N19").Value = -55.548037889039
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -55.555555555555
$ws.Range("I20").Value = 26
$ws.Range("J20").Value = 28
$ws.Range("K20").Value = -7.142857142857
$ws.Range("L20").Value = -33.333333333333
$ws.Range("M20").Value = 8.333333333333
$ws.Range("N20").Value = -94.104308390022
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 53
$ws.Range("E21").Value = -20.754716981132
$ws.Range("F21").Value = 135
$ws.Range("G21").Value = 173
$ws.Range("H21").Value = -21.965317919075
$ws.Range("I21").Value = 1041
$ws.Range("J21").Value = 1234
$ws.Range("K21").Value = -15.640194489465
$ws.Range("L21").Value = -24.235807860262
$ws.Range("M21").Value = 10.626992561105
$ws.Range("N21").Value = -66.698656429942
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -40
$ws.Range("I22").Value = 27
$ws.Range("J22").Value = 29
$ws.Range("K22").Value = -6.896551724137
$ws.Range("L22").Value = 3.846153846153
$ws.Range("M22").Value = -20.588235294117
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 48
$ws.Range("E24").Value = -20.833333333333
$ws.Range("F24").Value = 186
$ws.Range("G24").Value = 190
$ws.Range("H24").Value = -2.105263157894
$ws.Range("I24").Value = 1087
$ws.Range("J24").Value = 1285
$ws.Range("K24").Value = -15.408560311284
$ws.Range("L24").Value = -8.885163453478
$ws.Range("M24").Value = 13.821989528795
$ws.Range("C25").Value = 27
$ws.Range("D25").Value = 39
$ws.Range("E25").Value = -30.769230769230
$ws.Range("F25").Value = 135
$ws.Range("G25").Value = 159
$ws.Range("H25").Value = -15.094339622641
$ws.Range("I25").Value = 868
$ws.Range("J25").Value = 984
$ws.Range("K25").Value = -11.788617886178
$ws.Range("L25").Value = -9.583333333333
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 27
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = 17.391304347826
$ws.Range("I26").Value = 229
$ws.Range("J26").Value = 267
$ws.Range("K26").Value = -14.232209737827
$ws.Range("L26").Value = -5.371900826446
$ws.Range("M26").Value = 48.701298701298
$ws.Range("H27").Value = -100
$ws.Range("F28").Value = 8
$ws.Range("H28").Value = 300
$ws.Range("I28").Value = 46
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 15
$ws.Range("L28").Value = 4.545454545454

# --- Cells that change between numeric and text representation ---
# F27: was a number (1), becomes the literal text "0" (same display as the
#      workbook's "not applicable" placeholder used elsewhere in the sheet).
$ws.Range("F27").Value = "'0"
$ws.Range("C27").Copy()
$ws.Range("F27").PasteSpecial(-4122)

# D28/E28: were the "not applicable" placeholders, become real numbers.
$ws.Range("D28").Value = 1
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").Value = 200
$ws.Range("N15").Copy()
$ws.Range("E28").PasteSpecial(-4122)
